# Update cryptos list (prices in column D, 1h-volume % in column E).
# Values are stored as text (not numbers) to preserve exact source
# formatting (thousands-separator dots, trailing zeros, padding spaces,
# subscript digits, etc). Numeric-looking price strings are written
# with a leading apostrophe so Excel keeps them as text instead of
# silently coercing them to a Double (which would drop trailing zeros
# such as "4.40" -> "4.4").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.884.21"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = "'2.419.17"
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'551.19"
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = "'137.49"
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +3.21%  '
$ws.Range("D9").Value = "'0.105"
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("D10").Value = "'5.67"
$ws.Range("E10").Value = '  -3.22%  '
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").Value = "'25.44"
$ws.Range("E13").Value = '  +3.36%  '
$ws.Range("D14").Value = "'2.852.48"
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = "'59.839.43"
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = "'2.452.97"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = "'11.36"
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").Value = "'4.40"
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = "'329.35"
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = '  -5.00%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = "'66.48"
$ws.Range("E23").Value = '  +2.55%  '
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").Value = "'8.68"
$ws.Range("E25").Value = '  +2.94%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").Value = "'0.0₃0775"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").Value = "'1.78"
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").Value = "'169.31"
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = '  -1.96%  '
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("E38").Value = '  -2.62%  '
$ws.Range("D39").Value = "'0.409"
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("D40").Value = "'313.71"
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("D42").Value = "'138.60"
$ws.Range("E42").Value = '  -2.82%  '
$ws.Range("D43").Value = "'0.0966"
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("D44").Value = "'0.0519"
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = "'19.61"
$ws.Range("E45").Value = '  +2.35%  '
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").Value = "'0.0224"
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").Value = "'0.389"
$ws.Range("E48").Value = '  -4.82%  '
$ws.Range("D49").Value = "'17.68"
$ws.Range("E49").Value = '  -0.78%  '
$ws.Range("D50").Value = "'11.07"
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("E51").Value = '  -1.09%  '
